$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B2").Value = 0.6273408239700374
$summary.Range("C2").Value = 0.5739130434782609
$summary.Range("D2").Value = 0.9887640449438202
$summary.Range("E2").Value = 0.7262723521320495
$summary.Range("F2").Value = 0.8638743455497382
$summary.Range("G2").Value = 0.9620182200420463
$summary.Range("H2").Value = 0.7767783248467505
$summary.Range("I2").Value = 528
$summary.Range("J2").Value = 392
$summary.Range("K2").Value = 142
$summary.Range("L2").Value = 6

# --- Sheet: Classification Report ---
$classRep = $wb.Worksheets.Item("Classification Report")

# Row 2 - label "0"
$classRep.Range("B2").Value = 0.9594594594594594
$classRep.Range("C2").Value = 0.2659176029962547
$classRep.Range("D2").Value = 0.4164222873900293

# Row 3 - label "1"
$classRep.Range("B3").Value = 0.5739130434782609
$classRep.Range("C3").Value = 0.9887640449438202
$classRep.Range("D3").Value = 0.7262723521320495

# Row 4 - label "accuracy"
$classRep.Range("B4").Value = 0.6273408239700374
$classRep.Range("C4").Value = 0.6273408239700374
$classRep.Range("D4").Value = 0.6273408239700374
$classRep.Range("E4").Value = 0.6273408239700374

# Row 5 - label "macro avg"
$classRep.Range("B5").Value = 0.7666862514688602
$classRep.Range("C5").Value = 0.6273408239700374
$classRep.Range("D5").Value = 0.5713473197610395

# Row 6 - label "weighted avg"
$classRep.Range("B6").Value = 0.7666862514688602
$classRep.Range("C6").Value = 0.6273408239700374
$classRep.Range("D6").Value = 0.5713473197610395

# --- Sheet: Confusion Matrix ---
$confMat = $wb.Worksheets.Item("Confusion Matrix")

# Row 2 - Actual 0
$confMat.Range("B2").Value = 142
$confMat.Range("C2").Value = 392

# Row 3 - Actual 1
$confMat.Range("B3").Value = 6
$confMat.Range("C3").Value = 528
